# Updated symbol list (crypto price/volume refresh) on Fri Dec 16 02:34:51 UTC 2022
# with GitHub Actions. Refreshes the "Price" column (and a few rank/volume
# labels) with the latest scraped values, and rotates the TigerCash/BitKan/
# HotbitToken/NitroEx/LEO/BTSEToken/One block up one rank (One moves from
# rank 17 to rank 23 with its new price).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Columns in this sheet store every value (including price numbers) as
    # literal text. Force the Text number format first so that numeric-
    # looking strings ("261.42", "0.006219", ...) aren't silently coerced
    # into floating point numbers (which would also mangle trailing zeros
    # like "0.0003200" or "0.00005920").
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Price column refresh (rows 2-17) ---
Set-TextValue $ws.Range("D2")  "261.42"
Set-TextValue $ws.Range("D3")  "22.81"
Set-TextValue $ws.Range("D4")  "6.185"
Set-TextValue $ws.Range("D5")  "0.06094"
Set-TextValue $ws.Range("D6")  "6.744"
Set-TextValue $ws.Range("D7")  "3.454"
Set-TextValue $ws.Range("D8")  "1.365"
Set-TextValue $ws.Range("D9")  "0.7976"
Set-TextValue $ws.Range("D10") "0.1586"
Set-TextValue $ws.Range("D11") "0.08033"
Set-TextValue $ws.Range("D12") "0.03417"
Set-TextValue $ws.Range("D13") "0.03071"
Set-TextValue $ws.Range("D14") "0.09316"
Set-TextValue $ws.Range("D15") "3.861"
Set-TextValue $ws.Range("D16") "0.001693"
Set-TextValue $ws.Range("D17") "0.04837"

# --- Rows 18-24: ranking block rotates up by one (One wraps to the bottom) ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D18") "0.006219"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D19") "0.001092"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "0.003406"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D21") "0.0001499"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D22") "3.712"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D23") "2.247"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D24") "0.01337"
$ws.Range("E24").Value = "23OneONE"

# --- Remaining price refreshes further down the sheet ---
Set-TextValue $ws.Range("D25") "0.3364"
Set-TextValue $ws.Range("D27") "0.0003200"
Set-TextValue $ws.Range("D40") "0.04576"
Set-TextValue $ws.Range("D41") "0.007127"
Set-TextValue $ws.Range("D42") "0.003898"
Set-TextValue $ws.Range("D43") "0.1117"
Set-TextValue $ws.Range("D44") "0.01064"
Set-TextValue $ws.Range("D45") "0.002969"
Set-TextValue $ws.Range("D46") "0.00005920"
Set-TextValue $ws.Range("D48") "0.6996"
Set-TextValue $ws.Range("D49") "0.07644"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
Set-TextValue $ws.Range("D50") "0.00002099"
